$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# New experiment data: graph_id, size, no_comms, inside_prob, outside_prob
# (row, A, B, C, D, E) - rows 2..30 are rewritten in place (previously they
# referenced shared strings for D/E; now they hold plain numbers), and rows
# 31..34 are brand new experiment rows appended below the old data.
$data = @(
    @(2, 0, 100, 2, 0.4, 0.01),
    @(3, 1, 100, 2, 0.3, 0.02),
    @(4, 2, 100, 2, 0.5, 0.05),
    @(5, 3, 100, 2, 0.4, 0.08),
    @(6, 4, 100, 2, 0.4, 0.1),
    @(7, 5, 100, 2, 0.7, 0.12),
    @(8, 6, 100, 3, 0.35, 0.06),
    @(9, 7, 100, 3, 0.3, 0.02),
    @(10, 8, 100, 3, 0.4, 0.05),
    @(11, 9, 100, 3, 0.4, 0.02),
    @(12, 10, 100, 3, 0.4, 0.05),
    @(13, 11, 100, 3, 0.4, 0.08),
    @(14, 12, 100, 4, 0.4, 0.01),
    @(15, 13, 100, 4, 0.3, 0.03),
    @(16, 14, 100, 4, 0.5, 0.02),
    @(17, 15, 100, 4, 0.3, 0.02),
    @(18, 16, 100, 4, 0.3, 0.02),
    @(19, 17, 100, 4, 0.35, 0.02),
    @(20, 18, 100, 4, 0.5, 0.07),
    @(21, 19, 100, 5, 0.5, 0.05),
    @(22, 20, 100, 5, 0.4, 0.01),
    @(23, 21, 100, 5, 0.3, 0.02),
    @(24, 22, 100, 5, 0.3, 0.02),
    @(25, 23, 100, 5, 0.4, 0.05),
    @(26, 24, 100, 5, 0.5, 0.08),
    @(27, 25, 100, 6, 0.45, 0.1),
    @(28, 26, 100, 6, 0.4, 0.01),
    @(29, 27, 100, 6, 0.5, 0.02),
    @(30, 28, 100, 6, 0.4, 0.05),
    @(31, 29, 100, 6, 0.4, 0.03),
    @(32, 30, 100, 6, 0.8, 0.06),
    @(33, 31, 100, 6, 0.65, 0.08),
    @(34, 32, 100, 6, 0.6, 0.05)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# Restore the header row in case it was ever touched (defensive no-op).
$ws.Cells.Item(1, 1).Value = "graph_id"
$ws.Cells.Item(1, 2).Value = "size"
$ws.Cells.Item(1, 3).Value = "no_comms"
$ws.Cells.Item(1, 4).Value = "inside_prob"
$ws.Cells.Item(1, 5).Value = "outside_prob"

# View state: selection now covers the whole "size" column of the
# (larger) data range, and the window is zoomed to 100%.
$ws.Range("B2:B34").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
